$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# New translation rows appended at the bottom of the "Import" sheet (rows 38-41).
$ws.Range("A38").Value = "cs"
$ws.Range("A39").Value = "cs"
$ws.Range("B38").Value = "public.sign-out.title"
$ws.Range("B39").Value = "public.sign-out"
$ws.Range("C38").Value = "Odhlašování"
$ws.Range("C39").Value = "Probíhá odhlašování z aplikace, prosím vyčkejte…"

$ws.Range("A40").Value = "cs"
$ws.Range("B40").Value = "error.Who are you?"
$ws.Range("C40").Value = "Je nám líto, ale aplikace vás nepoznává."

$ws.Range("A41").Value = "cs"
$ws.Range("B41").Value = "error.Unknown login"
$ws.Range("C41").Value = "Přihlášení selhalo, zkontrolujte si prosím jméno a heslo."

# Carry the existing row style (row 37, the last pre-existing data row) down
# onto the four new rows, matching how the other data rows are formatted.
$ws.Range("A37:C37").Copy()
$ws.Range("A38:C41").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the updated view/selection state recorded in the sheet.
$ws.Range("B34").Select()
$excel.ActiveWindow.ScrollRow = 19
